# Apply corrected results after fixing error estimation and number of
# projected years (trends not yet rerun after filtering bug fix).

$wb = $excel.ActiveWorkbook

# --- Sheet: "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Range("C2").Value = 23
$ws1.Range("D2").Value = 2.4

$ws1.Range("C3").Value = 38
$ws1.Range("D3").Value = 7.3
$ws1.Range("E3").Value = 31.7

$ws1.Range("B4").Value = 21
$ws1.Range("C4").Value = 41
$ws1.Range("D4").Value = 51.2
$ws1.Range("E4").Value = 34.2

$ws1.Range("D5").Value = 12.2
$ws1.Range("E5").Value = 6.7

$ws1.Range("C6").Value = 10
$ws1.Range("D6").Value = 26.8
$ws1.Range("E6").Value = 8.300000000000001

$ws1.Range("B7").Value = 70
$ws1.Range("C7").Value = 134

# --- Sheet: "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")

$ws4.Range("C3").Value = 41
$ws4.Range("C4").Value = 120

# --- Sheet: "Interannual update - High Pri" ---
$ws5 = $wb.Worksheets.Item("Interannual update - High Pri")

$ws5.Range("B2").Value = 76
$ws5.Range("C2").Value = 73.8
$ws5.Range("D2").Value = 76
$ws5.Range("E2").Value = 88.40000000000001

# Insert a new row at position 3 ("Trend Different"), pushing the
# existing "IUCN" row down to row 4.
$ws5.Rows.Item(3).Insert()

$ws5.Range("A3").Value = "Trend Different"
$ws5.Range("B3").Value = 1
$ws5.Range("C3").Value = 1

$ws5.Range("A4").Value = "IUCN"
$ws5.Range("B4").Value = 26
$ws5.Range("C4").Value = 25.2
$ws5.Range("D4").Value = 10
$ws5.Range("E4").Value = 11.6
